$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.797.34'
$ws.Range("E2").Value = '  +1.14%  '
$ws.Range("D3").Value = '2.154.17'
$ws.Range("E3").Value = '  +2.29%  '
$ws.Range("E4").Value = '  +0.43%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '352.37'
$ws.Range("E5").Value = '  +5.44%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.009'
$ws.Range("E6").Value = '  +0.33%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5294'
$ws.Range("E7").Value = '  +1.31%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4586'
$ws.Range("E8").Value = '  +1.23%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '54.27'
$ws.Range("E9").Value = '  +1.70%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.09232'
$ws.Range("E10").Value = '  +3.44%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.189'
$ws.Range("E11").Value = '  +0.27%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '25.02'
$ws.Range("E12").Value = '  +3.74%  '
$ws.Range("D13").Value = '2.147.51'
$ws.Range("E13").Value = '  +2.16%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.928'
$ws.Range("E14").Value = '  +1.59%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.200'
$ws.Range("E15").Value = '  +2.13%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '102.46'
$ws.Range("E16").Value = '  +5.85%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001184'
$ws.Range("E17").Value = '  +3.62%  '
$ws.Range("E18").Value = '  +0.29%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06728'
$ws.Range("E19").Value = '  +1.08%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.67'
$ws.Range("E21").Value = '  +0.31%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.394'
$ws.Range("E22").Value = '  +0.86%  '
$ws.Range("D23").Value = '30.862.72'
$ws.Range("E23").Value = '  +1.10%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.93'
$ws.Range("E24").Value = '  +3.64%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.397'
$ws.Range("E25").Value = '  +1.73%  '
$ws.Range("D26").Value = '2.404.00'
$ws.Range("E26").Value = '  +2.29%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.71'
$ws.Range("E27").Value = '  +1.85%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.644'
$ws.Range("E28").Value = '  +4.62%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '165.20'
$ws.Range("E29").Value = '  +1.36%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '137.69'
$ws.Range("E30").Value = '  +2.84%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.230'
$ws.Range("E31").Value = '  +1.89%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1086'
$ws.Range("E32").Value = '  +1.40%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.684'
$ws.Range("E33").Value = '  +1.62%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.422'
$ws.Range("E34").Value = '  -0.17%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.016'
$ws.Range("E35").Value = '  +1.76%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.179'
$ws.Range("E36").Value = '  +6.48%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '10.49'
$ws.Range("E37").Value = '  +0.48%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02666'
$ws.Range("E38").Value = '  +3.03%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06944'
$ws.Range("E39").Value = '  +1.44%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2347'
$ws.Range("E40").Value = '  +2.25%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '12.74'
$ws.Range("E41").Value = '  +0.05%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6988'
$ws.Range("E42").Value = '  +1.59%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.279'
$ws.Range("E43").Value = '  +2.42%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.93'
$ws.Range("E44").Value = '  +5.75%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.371'
$ws.Range("E45").Value = '  +2.42%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6514'
$ws.Range("E46").Value = '  +2.30%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.761'
$ws.Range("E47").Value = '  +2.65%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.00000000370'
$ws.Range("E48").Value = '  +4.73%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.263'
$ws.Range("E49").Value = '  +1.14%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '83.71'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07335'
$ws.Range("E51").Value = '  +2.60%  '
